$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C8").Value = "bbb"
$ws.Range("C8").Select() | Out-Null
